$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Fix the buggy PADRC formulas in columns O:AH of 'DRC-PADRC'!2.
#    They previously extrapolated off a ratio of peak-demand values; they now
#    use the same TREND(...)-based formula already used in columns D:N.
# ---------------------------------------------------------------------------
$padrc = $wb.Worksheets.Item("DRC-PADRC")

$cols = @("O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")
foreach ($col in $cols) {
    $cell = $col + "2"
    $ref1 = $col + "1"
    $padrc.Range($cell).Formula = "=TREND(Calculations!`$A`$3:`$B`$3,Calculations!`$A`$2:`$B`$2,'DRC-PADRC'!$ref1)-'DRC-BDRC'!$cell"
}

# ---------------------------------------------------------------------------
# 2) Restore / update each sheet's selection (and the active sheet/tab) to
#    match the saved view state.
# ---------------------------------------------------------------------------

# Calculations: selection is now just B3
$calc = $wb.Worksheets.Item("Calculations")
$calc.Activate() | Out-Null
$calc.Range("B3").Select() | Out-Null

# DRC-BDRC: scrolled so column F is the leftmost visible column, selection B2:AH2
$bdrc = $wb.Worksheets.Item("DRC-BDRC")
$bdrc.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 6
$bdrc.Range("B2:AH2").Select() | Out-Null

# DRC-PADRC: becomes the active tab, selection N2:AH2
$padrc.Activate() | Out-Null
$padrc.Range("N2:AH2").Select() | Out-Null
